# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Both sheets contain the same data, and rows 5,6,7,8,20,21 get their F
# value bumped to the newer scrape result.

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 48
    6  = 355
    7  = 10898
    8  = 416
    20 = 1123
    21 = 53
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
